$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column H ("Importe") holds amount values stored as text. The original scrape
# produced Spanish-style formatting (thousands ".", decimal ","); normalize to
# plain decimal-dot text (e.g. "3.295,00" -> "3295.00") while keeping the cells
# as text so the fix round-trips through the shared-string table like the source.
$ws.Range("H2:H168").NumberFormat = "@"

$ws.Range("H2").Value = "3295.00"
$ws.Range("H3").Value = "2339.00"
$ws.Range("H4").Value = "2900.00"
$ws.Range("H5").Value = "50.00"
$ws.Range("H6").Value = "122000.00"
$ws.Range("H7").Value = "10703.55"
$ws.Range("H8").Value = "69814.35"
$ws.Range("H9").Value = "558.96"
$ws.Range("H10").Value = "279.43"
$ws.Range("H11").Value = "11114.34"
$ws.Range("H12").Value = "20.08"
$ws.Range("H13").Value = "41.16"
$ws.Range("H14").Value = "240.00"
$ws.Range("H15").Value = "28653.45"
$ws.Range("H16").Value = "3264.00"
$ws.Range("H17").Value = "431.24"
$ws.Range("H18").Value = "178.50"
$ws.Range("H19").Value = "14002.62"
$ws.Range("H20").Value = "17.80"
$ws.Range("H21").Value = "8048.00"
$ws.Range("H22").Value = "254.00"
$ws.Range("H23").Value = "22000.00"
$ws.Range("H24").Value = "37785.60"
$ws.Range("H25").Value = "2790.00"
$ws.Range("H26").Value = "13703.31"
$ws.Range("H27").Value = "6338.14"
$ws.Range("H28").Value = "354.00"
$ws.Range("H29").Value = "1487.00"
$ws.Range("H30").Value = "21.00"
$ws.Range("H31").Value = "5987.50"
$ws.Range("H32").Value = "1708.88"
$ws.Range("H33").Value = "88.00"
$ws.Range("H34").Value = "1230.00"
$ws.Range("H35").Value = "34859.80"
$ws.Range("H36").Value = "1380.00"
$ws.Range("H37").Value = "380.00"
$ws.Range("H38").Value = "460.00"
$ws.Range("H39").Value = "548.68"
$ws.Range("H40").Value = "3659.24"
$ws.Range("H41").Value = "21.00"
$ws.Range("H42").Value = "464.69"
$ws.Range("H43").Value = "2640.00"
$ws.Range("H44").Value = "754.20"
$ws.Range("H45").Value = "3021.54"
$ws.Range("H46").Value = "4527.18"
$ws.Range("H47").Value = "437.00"
$ws.Range("H48").Value = "226.60"
$ws.Range("H49").Value = "346.00"
$ws.Range("H50").Value = "2500.00"
$ws.Range("H51").Value = "375.00"
$ws.Range("H52").Value = "181.50"
$ws.Range("H53").Value = "1000.00"
$ws.Range("H54").Value = "19320.00"
$ws.Range("H55").Value = "328.00"
$ws.Range("H56").Value = "110.00"
$ws.Range("H57").Value = "1402.21"
$ws.Range("H58").Value = "720.00"
$ws.Range("H59").Value = "2218.00"
$ws.Range("H60").Value = "6165.00"
$ws.Range("H61").Value = "2780.00"
$ws.Range("H62").Value = "4000.00"
$ws.Range("H63").Value = "143043.05"
$ws.Range("H64").Value = "150000.00"
$ws.Range("H65").Value = "680822.45"
$ws.Range("H66").Value = "20267.50"
$ws.Range("H67").Value = "28.76"
$ws.Range("H68").Value = "1209.65"
$ws.Range("H69").Value = "1363.00"
$ws.Range("H70").Value = "400.00"
$ws.Range("H71").Value = "339.50"
$ws.Range("H72").Value = "1551.90"
$ws.Range("H73").Value = "20.80"
$ws.Range("H74").Value = "400.97"
$ws.Range("H75").Value = "35.00"
$ws.Range("H76").Value = "8.99"
$ws.Range("H77").Value = "18.00"
$ws.Range("H78").Value = "83289.50"
$ws.Range("H79").Value = "231.00"
$ws.Range("H80").Value = "110.00"
$ws.Range("H81").Value = "3895.00"
$ws.Range("H82").Value = "465.00"
$ws.Range("H83").Value = "507.77"
$ws.Range("H84").Value = "278.88"
$ws.Range("H85").Value = "37.40"
$ws.Range("H86").Value = "17357.00"
$ws.Range("H87").Value = "720.00"
$ws.Range("H88").Value = "2804.00"
$ws.Range("H89").Value = "129.00"
$ws.Range("H90").Value = "8302.00"
$ws.Range("H91").Value = "127.60"
$ws.Range("H92").Value = "500.00"
$ws.Range("H93").Value = "1920.00"
$ws.Range("H94").Value = "1858.00"
$ws.Range("H95").Value = "1400.00"
$ws.Range("H96").Value = "1333.36"
$ws.Range("H97").Value = "1535.00"
$ws.Range("H98").Value = "11384.80"
$ws.Range("H99").Value = "98250.00"
$ws.Range("H100").Value = "3000.00"
$ws.Range("H101").Value = "2000.00"
$ws.Range("H102").Value = "700.00"
$ws.Range("H103").Value = "2000.00"
$ws.Range("H104").Value = "768.00"
$ws.Range("H105").Value = "600.00"
$ws.Range("H106").Value = "1000.00"
$ws.Range("H107").Value = "4000.00"
$ws.Range("H108").Value = "13147.66"
$ws.Range("H109").Value = "965.00"
$ws.Range("H110").Value = "3000.00"
$ws.Range("H111").Value = "1900.00"
$ws.Range("H112").Value = "750.00"
$ws.Range("H113").Value = "4900.00"
$ws.Range("H114").Value = "1500.00"
$ws.Range("H115").Value = "400.00"
$ws.Range("H116").Value = "520.00"
$ws.Range("H117").Value = "4000.00"
$ws.Range("H118").Value = "30.00"
$ws.Range("H119").Value = "615.00"
$ws.Range("H120").Value = "570.00"
$ws.Range("H121").Value = "1270.00"
$ws.Range("H122").Value = "4820.00"
$ws.Range("H123").Value = "185.96"
$ws.Range("H124").Value = "79.23"
$ws.Range("H125").Value = "9210.00"
$ws.Range("H126").Value = "900.00"
$ws.Range("H127").Value = "716.80"
$ws.Range("H128").Value = "2830.00"
$ws.Range("H129").Value = "977.44"
$ws.Range("H130").Value = "27830.00"
$ws.Range("H131").Value = "2008.00"
$ws.Range("H132").Value = "1335.67"
$ws.Range("H133").Value = "284.00"
$ws.Range("H134").Value = "1572.35"
$ws.Range("H135").Value = "9413.00"
$ws.Range("H136").Value = "4318.11"
$ws.Range("H137").Value = "798.90"
$ws.Range("H138").Value = "1505.88"
$ws.Range("H139").Value = "402.00"
$ws.Range("H140").Value = "900.00"
$ws.Range("H141").Value = "341.00"
$ws.Range("H142").Value = "815.38"
$ws.Range("H143").Value = "3475.00"
$ws.Range("H144").Value = "300.22"
$ws.Range("H145").Value = "890.00"
$ws.Range("H146").Value = "56870.72"
$ws.Range("H147").Value = "18330.00"
$ws.Range("H148").Value = "55.81"
$ws.Range("H149").Value = "522697.37"
$ws.Range("H150").Value = "1650.00"
$ws.Range("H151").Value = "1625.00"
$ws.Range("H152").Value = "178900.00"
$ws.Range("H153").Value = "17500.00"
$ws.Range("H154").Value = "20000.00"
$ws.Range("H155").Value = "202500.00"
$ws.Range("H156").Value = "62500.00"
$ws.Range("H157").Value = "65000.00"
$ws.Range("H158").Value = "223000.00"
$ws.Range("H159").Value = "32000.00"
$ws.Range("H160").Value = "196300.00"
$ws.Range("H161").Value = "223000.00"
$ws.Range("H162").Value = "700.00"
$ws.Range("H163").Value = "29500.00"
$ws.Range("H164").Value = "12150.00"
$ws.Range("H165").Value = "2125.00"
$ws.Range("H166").Value = "544.00"
$ws.Range("H167").Value = "6000.00"
$ws.Range("H168").Value = "450.00"

# Fix stray commas used as separators between co-contractor names (should be
# periods); also drop dots from the "S.H." abbreviation in one entry.
$ws.Cells.Replace("IZAGUIRRE CARLOS MARIA, MOREND MARIA ELENA Y MOREND MARIA TERESA", "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA")
$ws.Cells.Replace("SCHAB DARIO, PEROTTI XAVIER, BENINCA MATIAS S.H.", "SCHAB DARIO. PEROTTI XAVIER. BENINCA MATIAS SH")
$ws.Cells.Replace("FERNANDEZ MARIO H, GALLICET OSCAR M", "FERNANDEZ MARIO H. GALLICET OSCAR M")
$ws.Cells.Replace("MARSICO GUILLERMO MIGUEL, MARSICO JUAN EDUARDO", "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO")

